# Updated cryptos list on Sun Oct 22 11:27:53 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2;  D = "29.939.28"; E = "  +0.62%  " },
    @{ Row = 3;  D = "1.633.44";  E = "  +1.85%  " },
    @{ Row = 4;  D = $null;       E = "  +0.17%  " },
    @{ Row = 5;  D = "214.47";    E = "  +0.92%  " },
    @{ Row = 6;  D = "0.518";     E = "  +0.02%  " },
    @{ Row = 7;  D = $null;       E = "  +0.16%  " },
    @{ Row = 8;  D = "28.48";     E = "  +0.04%  " },
    @{ Row = 9;  D = $null;       E = "  +1.14%  " },
    @{ Row = 10; D = "0.0608";    E = "  +0.67%  " },
    @{ Row = 11; D = "0.0908";    E = "  +0.31%  " },
    @{ Row = 12; D = "1.867.00";  E = "  +1.88%  " },
    @{ Row = 13; D = "1.636.38";  E = "  +1.85%  " },
    @{ Row = 14; D = "0.563";     E = "  +1.96%  " },
    @{ Row = 15; D = $null;       E = "  +15.22%  " },
    @{ Row = 16; D = "29.939.60"; E = "  +0.68%  " },
    @{ Row = 17; D = "3.85";      E = "  +1.79%  " },
    @{ Row = 18; D = "64.00";     E = "  -0.04%  " },
    @{ Row = 19; D = "241.88";    E = "  -0.10%  " },
    @{ Row = 20; D = "0.0₃0701";  E = "  +0.36%  " },
    @{ Row = 21; D = $null;       E = "  +0.13%  " },
    @{ Row = 22; D = $null;       E = "  +2.33%  " },
    @{ Row = 23; D = "9.80";      E = "  +3.49%  " },
    @{ Row = 24; D = $null;       E = "  +2.85%  " },
    @{ Row = 25; D = "158.84";    E = "  +2.38%  " },
    @{ Row = 26; D = "15.52";     E = "  +0.35%  " },
    @{ Row = 27; D = $null;       E = "  +0.44%  " },
    @{ Row = 28; D = "6.61";      E = "  +2.37%  " },
    @{ Row = 29; D = $null;       E = "  +0.17%  " },
    @{ Row = 30; D = $null;       E = "  +1.93%  " },
    @{ Row = 31; D = $null;       E = "  +3.95%  " },
    @{ Row = 32; D = $null;       E = "  +3.98%  " },
    @{ Row = 33; D = "3.18";      E = "  -0.25%  " },
    @{ Row = 34; D = "1.424.56";  E = $null },
    @{ Row = 35; D = "1.65";      E = "  +4.69%  " },
    @{ Row = 36; D = $null;       E = "  -0.83%  " },
    @{ Row = 38; D = $null;       E = "  -0.18%  " },
    @{ Row = 39; D = $null;       E = "  -0.11%  " },
    @{ Row = 40; D = "75.70";     E = "  +12.68%  " },
    @{ Row = 41; D = $null;       E = "  +1.18%  " },
    @{ Row = 42; D = $null;       E = "  +3.14%  " },
    @{ Row = 43; D = "0.828";     E = "  +1.29%  " },
    @{ Row = 44; D = "0.0492";    E = "  -0.48%  " },
    @{ Row = 45; D = $null;       E = "  +2.16%  " },
    @{ Row = 46; D = $null;       E = "  +0.28%  " },
    @{ Row = 47; D = "52.91";     E = "  -4.20%  " },
    @{ Row = 48; D = $null;       E = "  -0.55%  " },
    @{ Row = 49; D = "1.773.44";  E = "  +1.91%  " },
    @{ Row = 50; D = $null;       E = "  +9.51%  " },
    @{ Row = 51; D = "90.90";     E = "  +4.93%  " }
)

# Every Price cell in column D is stored as text (e.g. "214.38", "64.00"),
# including values that look numeric. Plain COM .Value assignment of a
# numeric-looking string gets auto-coerced to a real number (and can drop
# formatting like trailing zeros, e.g. "64.00" -> 64), so force the cell's
# number format to Text first to keep it consistent with the rest of the
# column.
foreach ($u in $updates) {
    $r = $u.Row
    if ($null -ne $u.D) {
        $cell = $ws.Range("D$r")
        $cell.NumberFormat = "@"
        $cell.Value = $u.D
    }
    if ($null -ne $u.E) {
        $ws.Range("E$r").Value = $u.E
    }
}
